# Cover Letter edits
$d = $word.ActiveDocument

# 1. Date change: "16.02." + "20" + "2" -> "17.02." + "202" (keep the trailing "3")
$d.Content.Find.Execute("16.02.20" + "2" + "3", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "17.02.2023", 2) | Out-Null

# 2. "Dear Prof" + ". " -> "Dear Prof. "
$d.Content.Find.Execute("Dear Prof. ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Dear Prof. ", 2) | Out-Null

# 3. closing quote + space, then "which we wish..." -> merge into one run
$d.Content.Find.Execute([char]8221 + " which we wish to be considered for publication in ", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, [char]8221 + " which we wish to be considered for publication in ", 2) | Out-Null

# 4. Abstract paragraph: remove "order to excite the field winding of " and split run
$d.Content.Find.Execute("more popular in order to excite the field winding of electrically", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "more popular in electrically", 2) | Out-Null

# 5. "...for publication. " + " " -> merge to "...for publication.  " (two trailing spaces)
$d.Content.Find.Execute("for publication.  ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "for publication.  ", 2) | Out-Null

Write-Output "done"
